$wb = $excel.ActiveWorkbook

# Overview sheet: "Latest HO Xliff Generate Date" for the cfbaec85... row (G4)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-09-05 02:50:56"

# zh-cn sheet: Correspond Handoff Datetime (H4) / Correspond Handback DateTime (K4)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-09-05 02:50:51"
$wsZhCn.Range("K4").Value = "2016-09-05 02:51:17"

# de-de sheet: Correspond Handoff Datetime (H4) matches the Overview value (shared string),
# Correspond Handback DateTime (K4) gets its own new value
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H4").Value = "2016-09-05 02:50:56"
$wsDeDe.Range("K4").Value = "2016-09-05 02:51:24"
